# Adapt column header formatting to respective input file names:
#   "<name>_old"  -> "<name>_FV2410"
#   "<name>_new"  -> "<name>_FV2504"
# Then turn the data range into an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the 21 header cells in row 1 (A1:U1).
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn A1:U66 into an Excel Table ("Table1") with the renamed headers,
#    so the table's column definitions (xl/tables/table1.xml) pick up the
#    new header text automatically.
$tableRange = $ws.Range("A1:U66")
$tbl = $ws.ListObjects.Add(1, $tableRange, [System.Type]::Missing, 1)
$tbl.Name = "Table1"

# 3) Freeze the header row (split above row 2 / below row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
